# [EXTRA SCRAPE] full data scraped for extra batting and bowling fields
#
# Adds a new worksheet "ODI Batting Extra" after the existing "ODI Batting"
# sheet, containing MATCH_CODE / BATTING_POSITION / NUM_4 / NUM_6 /
# PERCENT_RUNS_OF_TOTAL / MAN_OF_MATCH columns for two matches.

$wb = $excel.ActiveWorkbook

# Reference sheet whose header formatting (bold, bordered, centered) we reuse
$srcSheet = $wb.Worksheets.Item("ODI Batting")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Insert the new worksheet at the end of the workbook (after the last sheet)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ODI Batting Extra"

# Match the page margins used throughout the rest of the workbook
# (0.75in left/right, 1in top/bottom, 0.5in header/footer).
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# Copy the header row formatting (bold font + thin border + centered align)
# from the "ODI Batting" sheet's header row so the new header matches style.
$srcSheet.Range("A1:F1").Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Header row ---
$ws.Range("A1").Value = "MATCH_CODE"
$ws.Range("B1").Value = "BATTING_POSITION"
$ws.Range("C1").Value = "NUM_4"
$ws.Range("D1").Value = "NUM_6"
$ws.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$ws.Range("F1").Value = "MAN_OF_MATCH"

# --- Data rows ---
# Force a text number format for the text-valued data cells so values that
# look numeric (match codes, counts, percentages) are preserved as literal
# text, matching the scraped source data rather than being auto-converted
# by Excel. BATTING_POSITION (row 2, column B) is a genuine number, so that
# single cell is left out of this range and keeps the default "General"
# format.
$ws.Range("A2:A3").NumberFormat = "@"
$ws.Range("C2:F3").NumberFormat = "@"
$ws.Range("B3").NumberFormat = "@"

# Row 2: match 4735
$ws.Range("A2").Value = "4735"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "2"
$ws.Range("D2").Value = "0"
$ws.Range("E2").Value = "5.11%"
$ws.Range("F2").Value = "NO"

# Row 3: match 4745 (only MATCH_CODE and MAN_OF_MATCH populated)
$ws.Range("A3").Value = "4745"
$ws.Range("B3").Value = ""
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = "NO"
